$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.776.37'
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '2.357.24'
$ws.Range("E3").Value = '  -4.23%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.30'
$ws.Range("E5").Value = '  -1.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.50'
$ws.Range("E6").Value = '  -6.07%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -11.05%  '

$ws.Range("D9").Value = '2.356.51'
$ws.Range("E9").Value = '  -4.15%  '

$ws.Range("E10").Value = '  -1.87%  '

$ws.Range("E11").Value = '  +0.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("E12").Value = '  -3.56%  '

$ws.Range("E13").Value = '  -3.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.75'
$ws.Range("E14").Value = '  -5.18%  '

$ws.Range("D15").Value = '2.783.72'
$ws.Range("E15").Value = '  -4.08%  '

$ws.Range("D16").Value = '60.798.64'
$ws.Range("E16").Value = '  -0.63%  '

$ws.Range("E17").Value = '  -3.21%  '

$ws.Range("D18").Value = '2.358.23'
$ws.Range("E18").Value = '  -4.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.58'
$ws.Range("E19").Value = '  -4.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '315.16'
$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.57'
$ws.Range("E22").Value = '  -6.87%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.89'
$ws.Range("E24").Value = '  +2.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.15'
$ws.Range("E25").Value = '  -1.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.44'
$ws.Range("E26").Value = '  +10.86%  '

$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("D28").Value = '2.480.17'
$ws.Range("E28").Value = '  -3.80%  '

$ws.Range("D29").Value = '0.0₃0894'
$ws.Range("E29").Value = '  -6.94%  '

$ws.Range("E30").Value = '  -3.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '503.98'
$ws.Range("E31").Value = '  -8.19%  '

$ws.Range("E32").Value = '  -5.18%  '

$ws.Range("E33").Value = '  -1.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -5.93%  '

$ws.Range("E35").Value = '  -3.53%  '

$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.56'
$ws.Range("E37").Value = '  -5.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.51'
$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.371'
$ws.Range("E39").Value = '  -1.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.24'
$ws.Range("E40").Value = '  -10.56%  '

$ws.Range("E41").Value = '  -0.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.54'
$ws.Range("E42").Value = '  -1.35%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.14'
$ws.Range("E44").Value = '  -0.74%  '

$ws.Range("E45").Value = '  -7.74%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.51'
$ws.Range("E46").Value = '  -5.29%  '

$ws.Range("E47").Value = '  -2.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0510'
$ws.Range("E48").Value = '  -4.76%  '

$ws.Range("E49").Value = '  -9.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.569'
$ws.Range("E50").Value = '  -3.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0895'
$ws.Range("E51").Value = '  -4.26%  '
